$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: replace scenario content (Registro -> Capcha scenario) ---
$ws.Range("B4").Value = "Aparicion de capcha "
$ws.Range("C4").Value = "En la pagina de myShopify me intento registrar sin embargo el capcha no me deja debido a que es un robot"
$ws.Range("D4").Value = "Encontrarse en el formulario de registro."
$ws.Range("E4").Value = "1. ir a la opcion del formulario de registro. 2. llenar los datos entregados para el registro en el formulario. 3. darle submit. 4. verifico que si haya aprecido el capcha."
$ws.Range("F4").Value = "No registra mi usuario y aparece el capcha para verificar que no es un robot."

# --- Row 5: new scenario (Buscar funcional) ---
$ws.Range("B5").Value = "Buscar funcional"
$ws.Range("C5").Value = "En la pagina de myshopify en el home utilizo el buscador y me aparece un resultado relacionado"
$ws.Range("D5").Value = "Me encuentro en el home "
$ws.Range("E5").Value = "1. ingresar en el buscador  una palabra clave de algun articulo 2. ver el nombre de la primera opción y verificar que si contenga algo relaciónado con la palabra clave ingresada"
$ws.Range("F5").Value = "Me aparece un resultado parecido a lo que busque."

# --- Formatting: center align the whole working block (this also causes the
# previously "border only" style to become unused so it gets pruned, matching
# the saved workbook). Columns C, E, F wrap text; columns B, D do not. ---
$ws.Range("B5:F9").HorizontalAlignment = -4108
$ws.Range("B5:F9").VerticalAlignment = -4108

$ws.Range("C4:C5").WrapText = $true
$ws.Range("E4:E5").WrapText = $true
$ws.Range("F4:F5").WrapText = $true
$ws.Range("B4:B5").WrapText = $false
$ws.Range("D4:D5").WrapText = $false

# --- Row height for the new row ---
$ws.Rows.Item(5).RowHeight = 90

# --- Selection ends on B6, matching the post-edit cursor position ---
$ws.Range("B6").Select()
